$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @{ B = 0.8925;              C = 0.0025;  D = 0.865;              E = 0.9;     F = 0.8825; G = 0.8875;             H = 0.855; I = 0.8780000000000001 }
    3  = @{ B = 0.8475;              C = 0.005;   D = 0.835;              E = 0.855;   F = 0.8375; G = 0.8625;             H = 0.84;  I = 0.8459999999999999 }
    4  = @{ B = 0.825;               C = 0;       D = 0.835;              E = 0.825;   F = 0.8275; G = 0.83;               H = 0.79;  I = 0.8215 }
    5  = @{ B = 0.8275;              C = 0.0175;  D = 0.79;               E = 0.8125;  F = 0.8075; G = 0.8100000000000001; H = 0.77;  I = 0.798 }
    6  = @{ B = 0.82;                C = 0;       D = 0.8149999999999999; E = 0.7975;  F = 0.835;  G = 0.8325;             H = 0.8;   I = 0.8160000000000001 }
    7  = @{ B = 0.65;                C = 0.005;   D = 0.5649999999999999; E = 0.535;   F = 0.535;  G = 0.5375;             H = 0.515; I = 0.5375000000000001 }
    8  = @{ B = 0.9175;              C = 0.0025;  D = 0.885;              E = 0.89;    F = 0.89;   G = 0.8825;             H = 0.8575;I = 0.881 }
    9  = @{ B = 0.9275;              C = 0.01;    D = 0.89;               E = 0.9075;  F = 0.905;  G = 0.895;              H = 0.875; I = 0.8945000000000001 }
    10 = @{ B = 0.8925;              C = 0;       D = 0.8725000000000001; E = 0.87;    F = 0.88;   G = 0.885;              H = 0.85;  I = 0.8714999999999999 }
    11 = @{ B = 0.9275;              C = 0.01;    D = 0.895;              E = 0.905;   F = 0.905;  G = 0.8975;             H = 0.88;  I = 0.8965 }
}

$cols = @('B','C','D','E','F','G','H','I')

foreach ($row in $data.Keys) {
    $rowData = $data[$row]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $rowData[$col]
    }
}
